$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Clear the "Targetable" flags (column E) for rows 2-6; they previously held "y"
$ws.Range("E2:E6").ClearContents() | Out-Null

# Update the active cell/selection on the Parameters sheet to F11
$ws.Activate()
$ws.Range("F11").Select() | Out-Null
